$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC!row113
$ws.Range("H113").Value = 11945.286
$ws.Range("I113").Value = 17876.25
$ws.Range("J113").Value = 4037.3333
$ws.Range("K113").Value = 17876.25
$ws.Range("L113").Value = 4037.3333
$ws.Range("M113").Value = -14622.25
$ws.Range("N113").Value = -10545.3333

$ws = $wb.Worksheets.Item("ARM")
# ARM!row5
$ws.Range("H5").Value = 1111320.1
$ws.Range("I5").Value = 1250185.1
$ws.Range("J5").Value = 400
$ws.Range("K5").Value = 1250185.1
$ws.Range("L5").Value = 400
$ws.Range("M5").Value = -1250073.1
$ws.Range("N5").Value = -624

# ARM!row9
$ws.Range("H9").Value = 30000
$ws.Range("J9").Value = 20000
$ws.Range("L9").Value = 20000
$ws.Range("N9").Value = -20340

# ARM!row20
$ws.Range("H20").Value = 30000
$ws.Range("J20").Value = 20000
$ws.Range("L20").Value = 20000
$ws.Range("N20").Value = -20540

# ARM!row23
$ws.Range("H23").Value = 13900
$ws.Range("J23").Value = 13900
$ws.Range("L23").Value = 13900
$ws.Range("N23").Value = -14418

# ARM!row43
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

# ARM!row44
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

# ARM!row45
$ws.Range("H45").Value = 1500
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 1500
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 1500
$ws.Range("N45").Value = -2254
$ws.Range("M45").ClearContents()

# ARM!row55
$ws.Range("H55").Value = 70000
$ws.Range("J55").Value = 70000
$ws.Range("L55").Value = 70000
$ws.Range("N55").Value = -70630

# ARM!row61
$ws.Range("H61").Value = 2894.1516
$ws.Range("I61").Value = 2215
$ws.Range("J61").Value = 3533.353
$ws.Range("K61").Value = 2215
$ws.Range("L61").Value = 3533.353
$ws.Range("M61").Value = -2003
$ws.Range("N61").Value = -3957.353

# ARM!row97
$ws.Range("H97").Value = 8078.4614
$ws.Range("I97").Value = 8697.5
$ws.Range("K97").Value = 8697.5
$ws.Range("M97").Value = -8201.5

# ARM!row136
$ws.Range("H136").Value = 2894.1516
$ws.Range("I136").Value = 2215
$ws.Range("J136").Value = 3533.353
$ws.Range("K136").Value = 6645
$ws.Range("L136").Value = 10600.059
$ws.Range("M136").Value = -4095
$ws.Range("N136").Value = -15700.059

$ws = $wb.Worksheets.Item("BSM")
# BSM!row4
$ws.Range("H4").Value = 1111320.1
$ws.Range("I4").Value = 1250185.1
$ws.Range("J4").Value = 400
$ws.Range("K4").Value = 1250185.1
$ws.Range("L4").Value = 400
$ws.Range("M4").Value = -1250070.1
$ws.Range("N4").Value = -630

# BSM!row20
$ws.Range("H20").Value = 1360.0968
$ws.Range("I20").Value = 1146.6957
$ws.Range("J20").Value = 1973.625
$ws.Range("K20").Value = 1146.6957
$ws.Range("L20").Value = 1973.625
$ws.Range("M20").Value = -899.6957
$ws.Range("N20").Value = -2467.625

# BSM!row22
$ws.Range("H22").Value = 179.875
$ws.Range("I22").Value = 148.42857
$ws.Range("K22").Value = 148.42857
$ws.Range("M22").Value = 24.57142999999999

# BSM!row94
$ws.Range("H94").Value = 803.63635
$ws.Range("I94").Value = 885.55554
$ws.Range("J94").Value = 435
$ws.Range("K94").Value = 885.55554
$ws.Range("L94").Value = 435
$ws.Range("M94").Value = -434.55554
$ws.Range("N94").Value = -1337

# BSM!row107
$ws.Range("H107").Value = 839.2
$ws.Range("I107").Value = 798.5
$ws.Range("J107").Value = 866.3333
$ws.Range("K107").Value = 798.5
$ws.Range("L107").Value = 866.3333
$ws.Range("M107").Value = 1121.5
$ws.Range("N107").Value = -4706.3333

$ws = $wb.Worksheets.Item("CRP")
# CRP!row31
$ws.Range("H31").Value = 2017.2593
$ws.Range("I31").Value = 1247.8182
$ws.Range("J31").Value = 5402.8
$ws.Range("K31").Value = 1247.8182
$ws.Range("L31").Value = 5402.8
$ws.Range("M31").Value = -952.8181999999999
$ws.Range("N31").Value = -5992.8

# CRP!row34
$ws.Range("H34").Value = 2017.2593
$ws.Range("I34").Value = 1247.8182
$ws.Range("J34").Value = 5402.8
$ws.Range("K34").Value = 1247.8182
$ws.Range("L34").Value = 5402.8
$ws.Range("M34").Value = -1045.8182
$ws.Range("N34").Value = -5806.8

# CRP!row58
$ws.Range("H58").Value = 1406
$ws.Range("I58").Value = 722.5714
$ws.Range("J58").Value = 3798
$ws.Range("K58").Value = 722.5714
$ws.Range("L58").Value = 3798
$ws.Range("M58").Value = -519.5714
$ws.Range("N58").Value = -4204

# CRP!row136
$ws.Range("H136").Value = 1406
$ws.Range("I136").Value = 722.5714
$ws.Range("J136").Value = 3798
$ws.Range("K136").Value = 2167.7142
$ws.Range("L136").Value = 11394
$ws.Range("M136").Value = 382.2857999999997
$ws.Range("N136").Value = -16494

$ws = $wb.Worksheets.Item("CUL")
# CUL!row5
$ws.Range("H5").Value = 1528.6316
$ws.Range("I5").Value = 794.4167
$ws.Range("J5").Value = 2787.2856
$ws.Range("K5").Value = 2383.2501
$ws.Range("L5").Value = 8361.856800000001
$ws.Range("M5").Value = -2271.2501
$ws.Range("N5").Value = -8585.856800000001

# CUL!row113
$ws.Range("H113").Value = 27778850
$ws.Range("I113").Value = 399.2
$ws.Range("J113").Value = 38462870
$ws.Range("K113").Value = 1197.6
$ws.Range("L113").Value = 115388610
$ws.Range("M113").Value = 972.4000000000001
$ws.Range("N113").Value = -115392950

# CUL!row135
$ws.Range("H135").Value = 1528.6316
$ws.Range("I135").Value = 794.4167
$ws.Range("J135").Value = 2787.2856
$ws.Range("K135").Value = 7149.7503
$ws.Range("L135").Value = 25085.5704
$ws.Range("M135").Value = -4614.7503
$ws.Range("N135").Value = -30155.5704

$ws = $wb.Worksheets.Item("GSM")
# GSM!row43
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("N43").ClearContents()

# GSM!row46
$ws.Range("H46").Value = 13294
$ws.Range("I46").Value = 5238
$ws.Range("J46").Value = 21350
$ws.Range("K46").Value = 5238
$ws.Range("L46").Value = 21350
$ws.Range("M46").Value = -5082
$ws.Range("N46").Value = -21662

# GSM!row80
$ws.Range("H80").Value = 3021.375
$ws.Range("I80").Value = 2944.1667
$ws.Range("J80").Value = 3253
$ws.Range("K80").Value = 2944.1667
$ws.Range("L80").Value = 3253
$ws.Range("M80").Value = -1946.1667
$ws.Range("N80").Value = -5249

# GSM!row83
$ws.Range("H83").Value = 3021.375
$ws.Range("I83").Value = 2944.1667
$ws.Range("J83").Value = 3253
$ws.Range("K83").Value = 14720.8335
$ws.Range("L83").Value = 16265
$ws.Range("M83").Value = -9728.833500000001
$ws.Range("N83").Value = -26249

# GSM!row97
$ws.Range("H97").Value = 10750
$ws.Range("I97").Value = 1500
$ws.Range("J97").Value = 20000
$ws.Range("K97").Value = 1500
$ws.Range("L97").Value = 20000
$ws.Range("M97").Value = -1004
$ws.Range("N97").Value = -20992

$ws = $wb.Worksheets.Item("LTW")
# LTW!row68
$ws.Range("H68").Value = 1361.4445
$ws.Range("I68").Value = 1233.7333
$ws.Range("K68").Value = 1233.7333
$ws.Range("M68").Value = -484.7333000000001

# LTW!row71
$ws.Range("H71").Value = 1361.4445
$ws.Range("I71").Value = 1233.7333
$ws.Range("K71").Value = 6168.6665
$ws.Range("M71").Value = -2424.6665

# LTW!row82
$ws.Range("H82").Value = 972.4
$ws.Range("I82").Value = 553.3333
$ws.Range("J82").Value = 1601
$ws.Range("K82").Value = 553.3333
$ws.Range("L82").Value = 1601
$ws.Range("M82").Value = -192.3333
$ws.Range("N82").Value = -2323

# LTW!row85
$ws.Range("H85").Value = 972.4
$ws.Range("I85").Value = 553.3333
$ws.Range("J85").Value = 1601
$ws.Range("K85").Value = 553.3333
$ws.Range("L85").Value = 1601
$ws.Range("M85").Value = 694.6667
$ws.Range("N85").Value = -4097

# LTW!row93
$ws.Range("H93").Value = 664
$ws.Range("I93").Value = 674.9167
$ws.Range("J93").Value = 637.8
$ws.Range("K93").Value = 674.9167
$ws.Range("L93").Value = 637.8
$ws.Range("M93").Value = 573.0833
$ws.Range("N93").Value = -3133.8

$ws = $wb.Worksheets.Item("WVR")
# WVR!row126
$ws.Range("H126").Value = 50913.3
$ws.Range("I126").Value = 56237
$ws.Range("K126").Value = 168711
$ws.Range("M126").Value = -166241

# WVR!row138
$ws.Range("H138").Value = 80000
$ws.Range("J138").Value = 80000
$ws.Range("L138").Value = 80000
$ws.Range("N138").Value = -90280

# WVR!row139
$ws.Range("H139").Value = 50000
$ws.Range("J139").Value = 50000
$ws.Range("L139").Value = 50000
$ws.Range("N139").Value = -60280
